$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.607.76'
$ws.Range('E2').Value = '  -1.73%  '
$ws.Range('D3').Value = '1.666.04'
$ws.Range('E3').Value = '  -3.60%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.85'
$ws.Range('E5').Value = '  -1.94%  '
$ws.Range('E6').Value = '  -2.44%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.61'
$ws.Range('E8').Value = '  -2.20%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.262'
$ws.Range('E9').Value = '  -0.92%  '
$ws.Range('E10').Value = '  -1.76%  '
$ws.Range('E11').Value = '  -2.29%  '
$ws.Range('D12').Value = '1.901.07'
$ws.Range('E12').Value = '  -3.61%  '
$ws.Range('D13').Value = '1.674.32'
$ws.Range('E13').Value = '  -3.10%  '
$ws.Range('E14').Value = '  -3.31%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.562'
$ws.Range('E15').Value = '  -0.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.18'
$ws.Range('E16').Value = '  -2.05%  '
$ws.Range('D17').Value = '27.607.27'
$ws.Range('E17').Value = '  -1.60%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '243.14'
$ws.Range('E18').Value = '  -0.02%  '
$ws.Range('D19').Value = '0.0₃0728'
$ws.Range('E19').Value = '  -3.71%  '
$ws.Range('E20').Value = '  -4.30%  '
$ws.Range('E21').Value = '  -0.07%  '
$ws.Range('E22').Value = '  -3.56%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.29'
$ws.Range('E23').Value = '  -4.54%  '
$ws.Range('E24').Value = '  -4.61%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.71'
$ws.Range('E25').Value = '  -1.48%  '
$ws.Range('E26').Value = '  -4.37%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.42'
$ws.Range('E27').Value = '  -1.96%  '
$ws.Range('E28').Value = '  +0.04%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.111'
$ws.Range('E29').Value = '  -2.48%  '
$ws.Range('E30').Value = '  +2.84%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0500'
$ws.Range('E31').Value = '  -1.94%  '
$ws.Range('D33').Value = '1.470.01'
$ws.Range('E33').Value = '  -1.74%  '
$ws.Range('E34').Value = '  -5.12%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.56'
$ws.Range('E35').Value = '  -6.04%  '
$ws.Range('E36').Value = '  -1.64%  '
$ws.Range('E37').Value = '  -2.76%  '
$ws.Range('E38').Value = '  -1.51%  '
$ws.Range('E39').Value = '  -6.21%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '69.36'
$ws.Range('E40').Value = '  -1.97%  '
$ws.Range('E41').Value = '  -4.99%  '
$ws.Range('E42').Value = '  -0.05%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.41'
$ws.Range('E43').Value = '  -7.38%  '
$ws.Range('B44').Value = 'MXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.22'
$ws.Range('E44').Value = '  -3.84%  '
$ws.Range('D45').Value = '1.809.42'
$ws.Range('E45').Value = '  -3.53%  '
$ws.Range('E46').Value = '  -1.71%  '
$ws.Range('E47').Value = '  -4.50%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '89.30'
$ws.Range('E48').Value = '  -2.08%  '
$ws.Range('E49').Value = '  -4.49%  '
$ws.Range('E50').Value = '  -2.10%  '
$ws.Range('E51').Value = '  -4.53%  '
